$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
# Shape 1: Title 1
$tr = $s.Shapes.Item(1).TextFrame.TextRange
$tr.Characters(1, 18).Text = "{g0}=andray(){/g1}"
# Shape 2: Content Placeholder 2
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Characters(1, 23).Text = "{g0}ethay uickqay {/g1}"
$tr.Characters(24, 16).Text = "{g2}ownbray{/g3}"
$tr.Characters(40, 16).Text = "{g4} oxfay {/g5}"
$tr.Characters(56, 16).Text = "{g6}umpsjay{/g7}"
$tr.Characters(72, 24).Text = "{g8} overhay ethay {/g9}"
$tr.Characters(96, 17).Text = "{g10}azylay{/g11}"
$tr.Characters(113, 33).Text = "{g12} ogday. ethay uickqay {/g13}"
$tr.Characters(146, 18).Text = "{g14}ownbray{/g15}"
$tr.Characters(164, 18).Text = "{g16} oxfay {/g17}"
$tr.Characters(182, 18).Text = "{g18}umpsjay{/g19}"
$tr.Characters(200, 26).Text = "{g20} overhay ethay {/g21}"
$tr.Characters(226, 17).Text = "{g22}azylay{/g23}"
$tr.Characters(243, 33).Text = "{g24} ogday. ethay uickqay {/g25}"
$tr.Characters(276, 18).Text = "{g26}ownbray{/g27}"
$tr.Characters(294, 18).Text = "{g28} oxfay {/g29}"
$tr.Characters(312, 18).Text = "{g30}umpsjay{/g31}"
$tr.Characters(330, 26).Text = "{g32} overhay ethay {/g33}"
$tr.Characters(356, 17).Text = "{g34}azylay{/g35}"
$tr.Characters(373, 33).Text = "{g36} ogday. ethay uickqay {/g37}"
$tr.Characters(406, 18).Text = "{g38}ownbray{/g39}"
$tr.Characters(424, 18).Text = "{g40} oxfay {/g41}"
$tr.Characters(442, 18).Text = "{g42}umpsjay{/g43}"
$tr.Characters(460, 26).Text = "{g44} overhay ethay {/g45}"
$tr.Characters(486, 17).Text = "{g46}azylay{/g47}"
$tr.Characters(503, 33).Text = "{g48} ogday. ethay uickqay {/g49}"
$tr.Characters(536, 18).Text = "{g50}ownbray{/g51}"
$tr.Characters(554, 18).Text = "{g52} oxfay {/g53}"
$tr.Characters(572, 18).Text = "{g54}umpsjay{/g55}"
$tr.Characters(590, 26).Text = "{g56} overhay ethay {/g57}"
$tr.Characters(616, 17).Text = "{g58}azylay{/g59}"
$tr.Characters(633, 19).Text = "{g60} ogday. {/g61}"
$tr.Characters(654, 25).Text = "{g65}ethay uickqay {/g66}"
$tr.Characters(679, 18).Text = "{g67}ownbray{/g68}"
$tr.Characters(697, 18).Text = "{g69} oxfay {/g70}"
$tr.Characters(715, 18).Text = "{g71}umpsjay{/g72}"
$tr.Characters(733, 26).Text = "{g73} overhay ethay {/g74}"
$tr.Characters(759, 17).Text = "{g75}azylay{/g76}"
$tr.Characters(776, 33).Text = "{g77} ogday. ethay uickqay {/g78}"
$tr.Characters(809, 18).Text = "{g79}ownbray{/g80}"
$tr.Characters(827, 18).Text = "{g81} oxfay {/g82}"
$tr.Characters(845, 18).Text = "{g83}umpsjay{/g84}"
$tr.Characters(863, 26).Text = "{g85} overhay ethay {/g86}"
$tr.Characters(889, 17).Text = "{g87}azylay{/g88}"
$tr.Characters(906, 33).Text = "{g89} ogday. ethay uickqay {/g90}"
$tr.Characters(939, 18).Text = "{g91}ownbray{/g92}"
$tr.Characters(957, 18).Text = "{g93} oxfay {/g94}"
$tr.Characters(975, 18).Text = "{g95}umpsjay{/g96}"
$tr.Characters(993, 26).Text = "{g97} overhay ethay {/g98}"
$tr.Characters(1019, 18).Text = "{g99}azylay{/g100}"
$tr.Characters(1037, 35).Text = "{g101} ogday. ethay uickqay {/g102}"
$tr.Characters(1072, 20).Text = "{g103}ownbray{/g104}"
$tr.Characters(1092, 20).Text = "{g105} oxfay {/g106}"
$tr.Characters(1112, 20).Text = "{g107}umpsjay{/g108}"
$tr.Characters(1132, 28).Text = "{g109} overhay ethay {/g110}"
$tr.Characters(1160, 19).Text = "{g111}azylay{/g112}"
$tr.Characters(1179, 35).Text = "{g113} ogday. ethay uickqay {/g114}"
$tr.Characters(1214, 20).Text = "{g115}ownbray{/g116}"
$tr.Characters(1234, 20).Text = "{g117} oxfay {/g118}"
$tr.Characters(1254, 20).Text = "{g119}umpsjay{/g120}"
$tr.Characters(1274, 28).Text = "{g121} overhay ethay {/g122}"
$tr.Characters(1302, 19).Text = "{g123}azylay{/g124}"
$tr.Characters(1321, 21).Text = "{g125} ogday. {/g126}"
$tr.Characters(1344, 27).Text = "{g130}ethay uickqay {/g131}"
$tr.Characters(1371, 20).Text = "{g132}ownbray{/g133}"
$tr.Characters(1391, 20).Text = "{g134} oxfay {/g135}"
$tr.Characters(1411, 20).Text = "{g136}umpsjay{/g137}"
$tr.Characters(1431, 28).Text = "{g138} overhay ethay {/g139}"
$tr.Characters(1459, 19).Text = "{g140}azylay{/g141}"
$tr.Characters(1478, 35).Text = "{g142} ogday. ethay uickqay {/g143}"
$tr.Characters(1513, 20).Text = "{g144}ownbray{/g145}"
$tr.Characters(1533, 20).Text = "{g146} oxfay {/g147}"
$tr.Characters(1553, 20).Text = "{g148}umpsjay{/g149}"
$tr.Characters(1573, 28).Text = "{g150} overhay ethay {/g151}"
$tr.Characters(1601, 19).Text = "{g152}azylay{/g153}"
$tr.Characters(1620, 35).Text = "{g154} ogday. ethay uickqay {/g155}"
$tr.Characters(1655, 20).Text = "{g156}ownbray{/g157}"
$tr.Characters(1675, 20).Text = "{g158} oxfay {/g159}"
$tr.Characters(1695, 20).Text = "{g160}umpsjay{/g161}"
$tr.Characters(1715, 28).Text = "{g162} overhay ethay {/g163}"
$tr.Characters(1743, 19).Text = "{g164}azylay{/g165}"
$tr.Characters(1762, 35).Text = "{g166} ogday. ethay uickqay {/g167}"
$tr.Characters(1797, 20).Text = "{g168}ownbray{/g169}"
$tr.Characters(1817, 20).Text = "{g170} oxfay {/g171}"
$tr.Characters(1837, 20).Text = "{g172}umpsjay{/g173}"
$tr.Characters(1857, 28).Text = "{g174} overhay ethay {/g175}"
$tr.Characters(1885, 19).Text = "{g176}azylay{/g177}"
$tr.Characters(1904, 35).Text = "{g178} ogday. ethay uickqay {/g179}"
$tr.Characters(1939, 20).Text = "{g180}ownbray{/g181}"
$tr.Characters(1959, 20).Text = "{g182} oxfay {/g183}"
$tr.Characters(1979, 20).Text = "{g184}umpsjay{/g185}"
$tr.Characters(1999, 28).Text = "{g186} overhay ethay {/g187}"
$tr.Characters(2027, 19).Text = "{g188}azylay{/g189}"
$tr.Characters(2046, 21).Text = "{g190} ogday. {/g191}"
